$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

# --- LP1912: header updates ---
$ws1.Range("A2").Value = "Última actualización: 11:17:39"
$ws1.Range("A3").Value = "Total filas: 172"

# --- LP1912: row updates / new rows (90-91, 133-134, 139-177) ---
$ws1.Cells.Item(90,1).Value = "08:55:25"
$ws1.Cells.Item(90,2).Value = "09:16"
$ws1.Cells.Item(90,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(90,4).Value = 21
$ws1.Cells.Item(90,5).Value = "LP1912"
$ws1.Cells.Item(91,1).Value = "07:24:45"
$ws1.Cells.Item(91,2).Value = "09:16"
$ws1.Cells.Item(91,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(91,4).Value = 112
$ws1.Cells.Item(91,5).Value = "LP1912"
$ws1.Cells.Item(133,1).Value = "10:52:37"
$ws1.Cells.Item(133,2).Value = "11:06"
$ws1.Cells.Item(133,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(133,4).Value = 14
$ws1.Cells.Item(133,5).Value = "LP1912"
$ws1.Cells.Item(134,1).Value = "09:26:30"
$ws1.Cells.Item(134,2).Value = "11:06"
$ws1.Cells.Item(134,3).Value = "16_P MOR-167 Y 521"
$ws1.Cells.Item(134,4).Value = 100
$ws1.Cells.Item(134,5).Value = "LP1912"
$ws1.Cells.Item(139,1).Value = "11:17:39"
$ws1.Cells.Item(139,2).Value = "11:18"
$ws1.Cells.Item(139,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(139,4).Value = 1
$ws1.Cells.Item(139,5).Value = "LP1912"
$ws1.Cells.Item(140,1).Value = "11:17:39"
$ws1.Cells.Item(140,2).Value = "11:18"
$ws1.Cells.Item(140,3).Value = "15_ABASTO"
$ws1.Cells.Item(140,4).Value = 1
$ws1.Cells.Item(140,5).Value = "LP1912"
$ws1.Cells.Item(141,1).Value = "09:26:30"
$ws1.Cells.Item(141,2).Value = "11:19"
$ws1.Cells.Item(141,3).Value = "86_EST CHICA-ESC AGRARIA"
$ws1.Cells.Item(141,4).Value = 113
$ws1.Cells.Item(141,5).Value = "LP1912"
$ws1.Cells.Item(142,1).Value = "11:17:39"
$ws1.Cells.Item(142,2).Value = "11:20"
$ws1.Cells.Item(142,3).Value = "225_C ROCA-H SUR"
$ws1.Cells.Item(142,4).Value = 3
$ws1.Cells.Item(142,5).Value = "LP1912"
$ws1.Cells.Item(143,1).Value = "10:13:53"
$ws1.Cells.Item(143,2).Value = "11:20"
$ws1.Cells.Item(143,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(143,4).Value = 67
$ws1.Cells.Item(143,5).Value = "LP1912"
$ws1.Cells.Item(144,1).Value = "09:26:30"
$ws1.Cells.Item(144,2).Value = "11:21"
$ws1.Cells.Item(144,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(144,4).Value = 115
$ws1.Cells.Item(144,5).Value = "LP1912"
$ws1.Cells.Item(145,1).Value = "10:13:53"
$ws1.Cells.Item(145,2).Value = "11:26"
$ws1.Cells.Item(145,3).Value = "225_C ROCA-H SUR"
$ws1.Cells.Item(145,4).Value = 73
$ws1.Cells.Item(145,5).Value = "LP1912"
$ws1.Cells.Item(146,1).Value = "10:52:37"
$ws1.Cells.Item(146,2).Value = "11:27"
$ws1.Cells.Item(146,3).Value = "225_C ROCA-H SUR"
$ws1.Cells.Item(146,4).Value = 35
$ws1.Cells.Item(146,5).Value = "LP1912"
$ws1.Cells.Item(147,1).Value = "10:13:53"
$ws1.Cells.Item(147,2).Value = "11:32"
$ws1.Cells.Item(147,3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(147,4).Value = 79
$ws1.Cells.Item(147,5).Value = "LP1912"
$ws1.Cells.Item(148,1).Value = "10:52:37"
$ws1.Cells.Item(148,2).Value = "11:34"
$ws1.Cells.Item(148,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(148,4).Value = 42
$ws1.Cells.Item(148,5).Value = "LP1912"
$ws1.Cells.Item(149,1).Value = "10:52:37"
$ws1.Cells.Item(149,2).Value = "11:35"
$ws1.Cells.Item(149,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(149,4).Value = 43
$ws1.Cells.Item(149,5).Value = "LP1912"
$ws1.Cells.Item(150,1).Value = "10:13:53"
$ws1.Cells.Item(150,2).Value = "11:38"
$ws1.Cells.Item(150,3).Value = "10_OLMOS"
$ws1.Cells.Item(150,4).Value = 85
$ws1.Cells.Item(150,5).Value = "LP1912"
$ws1.Cells.Item(151,1).Value = "10:13:53"
$ws1.Cells.Item(151,2).Value = "11:41"
$ws1.Cells.Item(151,3).Value = "17_ROMERO"
$ws1.Cells.Item(151,4).Value = 88
$ws1.Cells.Item(151,5).Value = "LP1912"
$ws1.Cells.Item(152,1).Value = "10:52:37"
$ws1.Cells.Item(152,2).Value = "11:42"
$ws1.Cells.Item(152,3).Value = "17_ROMERO"
$ws1.Cells.Item(152,4).Value = 50
$ws1.Cells.Item(152,5).Value = "LP1912"
$ws1.Cells.Item(153,1).Value = "10:52:37"
$ws1.Cells.Item(153,2).Value = "11:43"
$ws1.Cells.Item(153,3).Value = "10_OLMOS"
$ws1.Cells.Item(153,4).Value = 51
$ws1.Cells.Item(153,5).Value = "LP1912"
$ws1.Cells.Item(154,1).Value = "10:13:53"
$ws1.Cells.Item(154,2).Value = "11:51"
$ws1.Cells.Item(154,3).Value = "215B_EL PATO"
$ws1.Cells.Item(154,4).Value = 98
$ws1.Cells.Item(154,5).Value = "LP1912"
$ws1.Cells.Item(155,1).Value = "11:17:39"
$ws1.Cells.Item(155,2).Value = "11:52"
$ws1.Cells.Item(155,3).Value = "15_ABASTO"
$ws1.Cells.Item(155,4).Value = 35
$ws1.Cells.Item(155,5).Value = "LP1912"
$ws1.Cells.Item(156,1).Value = "10:13:53"
$ws1.Cells.Item(156,2).Value = "11:58"
$ws1.Cells.Item(156,3).Value = "225_GOMEZ"
$ws1.Cells.Item(156,4).Value = 105
$ws1.Cells.Item(156,5).Value = "LP1912"
$ws1.Cells.Item(157,1).Value = "10:52:37"
$ws1.Cells.Item(157,2).Value = "11:59"
$ws1.Cells.Item(157,3).Value = "225_GOMEZ"
$ws1.Cells.Item(157,4).Value = 67
$ws1.Cells.Item(157,5).Value = "LP1912"
$ws1.Cells.Item(158,1).Value = "10:13:53"
$ws1.Cells.Item(158,2).Value = "12:02"
$ws1.Cells.Item(158,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(158,4).Value = 109
$ws1.Cells.Item(158,5).Value = "LP1912"
$ws1.Cells.Item(159,1).Value = "11:17:39"
$ws1.Cells.Item(159,2).Value = "12:04"
$ws1.Cells.Item(159,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(159,4).Value = 47
$ws1.Cells.Item(159,5).Value = "LP1912"
$ws1.Cells.Item(160,1).Value = "10:52:37"
$ws1.Cells.Item(160,2).Value = "12:06"
$ws1.Cells.Item(160,3).Value = "10_OLMOS"
$ws1.Cells.Item(160,4).Value = 74
$ws1.Cells.Item(160,5).Value = "LP1912"
$ws1.Cells.Item(161,1).Value = "10:13:53"
$ws1.Cells.Item(161,2).Value = "12:06"
$ws1.Cells.Item(161,3).Value = "14_ABASTO"
$ws1.Cells.Item(161,4).Value = 113
$ws1.Cells.Item(161,5).Value = "LP1912"
$ws1.Cells.Item(162,1).Value = "10:13:53"
$ws1.Cells.Item(162,2).Value = "12:06"
$ws1.Cells.Item(162,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(162,4).Value = 113
$ws1.Cells.Item(162,5).Value = "LP1912"
$ws1.Cells.Item(163,1).Value = "11:17:39"
$ws1.Cells.Item(163,2).Value = "12:13"
$ws1.Cells.Item(163,3).Value = "10_OLMOS"
$ws1.Cells.Item(163,4).Value = 56
$ws1.Cells.Item(163,5).Value = "LP1912"
$ws1.Cells.Item(164,1).Value = "10:52:37"
$ws1.Cells.Item(164,2).Value = "12:16"
$ws1.Cells.Item(164,3).Value = "14_ABASTO"
$ws1.Cells.Item(164,4).Value = 84
$ws1.Cells.Item(164,5).Value = "LP1912"
$ws1.Cells.Item(165,1).Value = "10:52:37"
$ws1.Cells.Item(165,2).Value = "12:20"
$ws1.Cells.Item(165,3).Value = "215A_EL PATO"
$ws1.Cells.Item(165,4).Value = 88
$ws1.Cells.Item(165,5).Value = "LP1912"
$ws1.Cells.Item(166,1).Value = "11:17:39"
$ws1.Cells.Item(166,2).Value = "12:20"
$ws1.Cells.Item(166,3).Value = "14_ABASTO"
$ws1.Cells.Item(166,4).Value = 63
$ws1.Cells.Item(166,5).Value = "LP1912"
$ws1.Cells.Item(167,1).Value = "10:52:37"
$ws1.Cells.Item(167,2).Value = "12:21"
$ws1.Cells.Item(167,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(167,4).Value = 89
$ws1.Cells.Item(167,5).Value = "LP1912"
$ws1.Cells.Item(168,1).Value = "10:52:37"
$ws1.Cells.Item(168,2).Value = "12:23"
$ws1.Cells.Item(168,3).Value = "17_ROMERO"
$ws1.Cells.Item(168,4).Value = 91
$ws1.Cells.Item(168,5).Value = "LP1912"
$ws1.Cells.Item(169,1).Value = "11:17:39"
$ws1.Cells.Item(169,2).Value = "12:34"
$ws1.Cells.Item(169,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(169,4).Value = 77
$ws1.Cells.Item(169,5).Value = "LP1912"
$ws1.Cells.Item(170,1).Value = "10:52:37"
$ws1.Cells.Item(170,2).Value = "12:36"
$ws1.Cells.Item(170,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(170,4).Value = 104
$ws1.Cells.Item(170,5).Value = "LP1912"
$ws1.Cells.Item(171,1).Value = "10:52:37"
$ws1.Cells.Item(171,2).Value = "12:38"
$ws1.Cells.Item(171,3).Value = "17_179 Y 38"
$ws1.Cells.Item(171,4).Value = 106
$ws1.Cells.Item(171,5).Value = "LP1912"
$ws1.Cells.Item(172,1).Value = "11:17:39"
$ws1.Cells.Item(172,2).Value = "12:41"
$ws1.Cells.Item(172,3).Value = "10_OLMOS"
$ws1.Cells.Item(172,4).Value = 84
$ws1.Cells.Item(172,5).Value = "LP1912"
$ws1.Cells.Item(173,1).Value = "11:17:39"
$ws1.Cells.Item(173,2).Value = "12:48"
$ws1.Cells.Item(173,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(173,4).Value = 91
$ws1.Cells.Item(173,5).Value = "LP1912"
$ws1.Cells.Item(174,1).Value = "11:17:39"
$ws1.Cells.Item(174,2).Value = "12:49"
$ws1.Cells.Item(174,3).Value = "17_ROMERO"
$ws1.Cells.Item(174,4).Value = 92
$ws1.Cells.Item(174,5).Value = "LP1912"
$ws1.Cells.Item(175,1).Value = "10:52:37"
$ws1.Cells.Item(175,2).Value = "12:50"
$ws1.Cells.Item(175,3).Value = "15_ABASTO"
$ws1.Cells.Item(175,4).Value = 118
$ws1.Cells.Item(175,5).Value = "LP1912"
$ws1.Cells.Item(176,1).Value = "11:17:39"
$ws1.Cells.Item(176,2).Value = "13:06"
$ws1.Cells.Item(176,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(176,4).Value = 109
$ws1.Cells.Item(176,5).Value = "LP1912"
$ws1.Cells.Item(177,1).Value = "11:17:39"
$ws1.Cells.Item(177,2).Value = "13:13"
$ws1.Cells.Item(177,3).Value = "215D_EL PATO"
$ws1.Cells.Item(177,4).Value = 116
$ws1.Cells.Item(177,5).Value = "LP1912"

# --- LP1912-215: header updates ---
$ws2.Range("A2").Value = "Última actualización: 11:17:39"
$ws2.Range("A3").Value = "Total filas: 22"

# --- LP1912-215: new row 27 ---
$ws2.Cells.Item(27,1).Value = "11:17:39"
$ws2.Cells.Item(27,2).Value = "13:13"
$ws2.Cells.Item(27,3).Value = "215D_EL PATO"
$ws2.Cells.Item(27,4).Value = 116
$ws2.Cells.Item(27,5).Value = "LP1912"

# --- 6203-6173: header updates ---
$ws3.Range("A2").Value = "Última actualización: 11:17:39"
$ws3.Range("A3").Value = "Total filas: 30"

# --- 6203-6173: new row 35 ---
$ws3.Cells.Item(35,1).Value = "11:17:39"
$ws3.Cells.Item(35,2).Value = "12:53"
$ws3.Cells.Item(35,3).Value = "215C_LA PLATA"
$ws3.Cells.Item(35,4).Value = 96
$ws3.Cells.Item(35,5).Value = "L6203"
